$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D entirely (no longer used in the updated table)
$ws.Columns.Item(4).Delete()

# Rewrite the full author/year/metric table with corrected and additional rows
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "Performance Metric"

$ws.Range("A2").Value = 2004
$ws.Range("B2").Value = "Geurs and van Wee"
$ws.Range("C2").Value = "Accessibility (isochrone, gravity, logsum)"

$ws.Range("A3").Value = 2006
$ws.Range("B3").Value = "Scott et al."
$ws.Range("C3").Value = "Travel time and cost"

$ws.Range("A4").Value = 2007
$ws.Range("B4").Value = "Abdel-Rahim et al."
$ws.Range("C4").Value = "Network Connectivity"

$ws.Range("A5").Value = 2008
$ws.Range("B5").Value = "Taylor, M"
$ws.Range("C5").Value = "Accessibility (logsum)"

$ws.Range("A6").Value = 2010
$ws.Range("B6").Value = "Peeta et al."
$ws.Range("C6").Value = "Travel time and cost"

$ws.Range("A7").Value = 2010
$ws.Range("B7").Value = "Geurs et al."
$ws.Range("C7").Value = "Accessibility (logsum)"

$ws.Range("A8").Value = 2010
$ws.Range("B8").Value = "Levinson and Zhu"
$ws.Range("C8").Value = "Travel time and cost"

$ws.Range("A9").Value = 2010
$ws.Range("B9").Value = "Zhu et al."
$ws.Range("C9").Value = "Travel time and cost"

$ws.Range("A10").Value = 2011
$ws.Range("B10").Value = "Agarwal et al."
$ws.Range("C10").Value = "Network connectivity"

$ws.Range("A11").Value = 2011
$ws.Range("B11").Value = "Ip and Wang"
$ws.Range("C11").Value = "Network connectivity"

$ws.Range("A12").Value = 2011
$ws.Range("B12").Value = "Serulle et al."
$ws.Range("C12").Value = "Travel time and cost"

$ws.Range("A13").Value = 2011
$ws.Range("B13").Value = "Ibrahim, S"
$ws.Range("C13").Value = "Travel time and cost"

$ws.Range("A14").Value = 2011
$ws.Range("B14").Value = "Xie and Levinson"
$ws.Range("C14").Value = "Accessibility (isochrone)"

$ws.Range("A15").Value = 2012
$ws.Range("B15").Value = "Jenelius and Mattson"
$ws.Range("C15").Value = "Travel time and cost"

$ws.Range("A16").Value = 2012
$ws.Range("B16").Value = "Taylor and Susilawati"
$ws.Range("C16").Value = "Accessibility (gravity)"

$ws.Range("A17").Value = 2013
$ws.Range("B17").Value = "Omer et al."
$ws.Range("C17").Value = "Travel time and cost"

$ws.Range("A18").Value = 2014
$ws.Range("B18").Value = "Balijepalli and Oppong"
$ws.Range("C18").Value = "Travel time and cost"

$ws.Range("A19").Value = 2014
$ws.Range("B19").Value = "Osei-Asamoah and Lownes"
$ws.Range("C19").Value = "Network connectivity"

$ws.Range("A20").Value = 2014
$ws.Range("B20").Value = "Guze"
$ws.Range("C20").Value = "Network connectivity"

$ws.Range("A21").Value = 2015
$ws.Range("B21").Value = "Zhang et al."
$ws.Range("C21").Value = "Network connectivity"

$ws.Range("A22").Value = 2015
$ws.Range("B22").Value = "Jaller et al."
$ws.Range("C22").Value = "Travel time and cost"

$ws.Range("A23").Value = 2015
$ws.Range("B23").Value = "Xu et al."
$ws.Range("C23").Value = "Network connectivity"

$ws.Range("A24").Value = 2016
$ws.Range("B24").Value = "Winkler, C."
$ws.Range("C24").Value = "Accessibility (gravity)"

$ws.Range("A25").Value = 2017
$ws.Range("B25").Value = "Ganin et al."
$ws.Range("C25").Value = "Accessibility (gravity)"

$ws.Range("A26").Value = 2019
$ws.Range("B26").Value = "Vodak et al."
$ws.Range("C26").Value = "Network connectivity"

$ws.Range("A27").Value = 2019
$ws.Range("B27").Value = "Hackl and Adey"
$ws.Range("C27").Value = "Network connectivity"

$ws.Range("A28").Value = 2019
$ws.Range("B28").Value = "Gecchele et al."
$ws.Range("C28").Value = "Accessibility (logsum)"

# Ensure all table cells (including the newly added rows) keep the
# original formatting (font + centered alignment) without introducing
# new style/font definitions - copy formats only from the header cell
$ws.Range("A1").Copy()
$ws.Range("A1:C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the selection to match the saved view state
$ws.Range("B19").Select() | Out-Null
